$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 287653.44
$ws.Range("J17").Value = 287653.44
$ws.Range("L17").Value = 862960.3200000001
$ws.Range("N17").Value = -863296.3200000001

# Row 38
$ws.Range("H38").Value = 1958.9524
$ws.Range("I38").Value = 466.18182
$ws.Range("K38").Value = 1398.54546
$ws.Range("M38").Value = -1026.54546

# Row 98
$ws.Range("H98").Value = 1597
$ws.Range("I98").Value = 1658.3077
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 1658.3077
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = -160.3077000000001
$ws.Range("N98").Value = -3796

# Row 107
$ws.Range("H107").Value = 362.06668
$ws.Range("I107").Value = 421.36365
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 421.36365
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1498.63635
$ws.Range("N107").Value = -4039

# Row 113
$ws.Range("H113").Value = 4264
$ws.Range("I113").Value = 3466.6667
$ws.Range("K113").Value = 3466.6667
$ws.Range("M113").Value = -212.6667000000002

# Row 122
$ws.Range("H122").Value = 1597
$ws.Range("I122").Value = 1658.3077
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 4974.9231
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -2524.9231
$ws.Range("N122").Value = -7300

# Row 127
$ws.Range("H127").Value = 1130.7142
$ws.Range("I127").Value = 1130.7142
$ws.Range("K127").Value = 3392.1426
$ws.Range("M127").Value = 1567.8574

# Row 137
$ws.Range("H137").Value = 2010.5306
$ws.Range("I137").Value = 1412.8889
$ws.Range("K137").Value = 4238.6667
$ws.Range("M137").Value = -1688.6667

# Row 138
$ws.Range("H138").Value = 3649.0986
$ws.Range("I138").Value = 3345.125
$ws.Range("J138").Value = 3737.5273
$ws.Range("K138").Value = 10035.375
$ws.Range("L138").Value = 11212.5819
$ws.Range("M138").Value = -4895.375
$ws.Range("N138").Value = -21492.5819

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 33665
$ws.Range("I6").Value = 21002
$ws.Range("J6").Value = 39996.5
$ws.Range("K6").Value = 21002
$ws.Range("L6").Value = 39996.5
$ws.Range("M6").Value = -20829
$ws.Range("N6").Value = -40342.5

# Row 32
$ws.Range("H32").Value = 9201.064
$ws.Range("I32").Value = 5798.8125
$ws.Range("J32").Value = 16709.482
$ws.Range("K32").Value = 5798.8125
$ws.Range("L32").Value = 16709.482
$ws.Range("M32").Value = -5511.8125
$ws.Range("N32").Value = -17283.482

# Row 61
$ws.Range("H61").Value = 55559130
$ws.Range("I61").Value = 57146230
$ws.Range("K61").Value = 57146230
$ws.Range("M61").Value = -57146018

# Row 74
$ws.Range("H74").Value = 34488090
$ws.Range("J74").Value = 2430
$ws.Range("L74").Value = 2430
$ws.Range("N74").Value = -4178

# Row 77
$ws.Range("H77").Value = 34488090
$ws.Range("J77").Value = 2430
$ws.Range("L77").Value = 12150
$ws.Range("N77").Value = -20886

# Row 88
$ws.Range("H88").Value = 2774.5557
$ws.Range("I88").Value = 2589.6365
$ws.Range("J88").Value = 3065.1428
$ws.Range("K88").Value = 2589.6365
$ws.Range("L88").Value = 3065.1428
$ws.Range("M88").Value = -2183.6365
$ws.Range("N88").Value = -3877.1428

# Row 91
$ws.Range("H91").Value = 2774.5557
$ws.Range("I91").Value = 2589.6365
$ws.Range("J91").Value = 3065.1428
$ws.Range("K91").Value = 2589.6365
$ws.Range("L91").Value = 3065.1428
$ws.Range("M91").Value = -1185.6365
$ws.Range("N91").Value = -5873.1428

# Row 110
$ws.Range("H110").Value = 86299
$ws.Range("J110").Value = 4337.6665
$ws.Range("L110").Value = 4337.6665
$ws.Range("N110").Value = -8427.666499999999

# Row 136
$ws.Range("H136").Value = 55559130
$ws.Range("I136").Value = 57146230
$ws.Range("K136").Value = 171438690
$ws.Range("M136").Value = -171436140

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1708.8572
$ws.Range("I20").Value = 1578.9524
$ws.Range("K20").Value = 1578.9524
$ws.Range("M20").Value = -1331.9524

# Row 105
$ws.Range("H105").Value = 4570
$ws.Range("I105").Value = 5105
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 5105
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -3358
$ws.Range("N105").Value = -6994

# Row 134
$ws.Range("H134").Value = 11907447
$ws.Range("I134").Value = 13160439
$ws.Range("J134").Value = 4028.25
$ws.Range("K134").Value = 39481317
$ws.Range("L134").Value = 12084.75
$ws.Range("M134").Value = -39478782
$ws.Range("N134").Value = -17154.75

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 127
$ws.Range("I2").Value = 86
$ws.Range("K2").Value = 86
$ws.Range("M2").Value = 27

# Row 31
$ws.Range("H31").Value = 5550.9165
$ws.Range("I31").Value = 4106.5
$ws.Range("K31").Value = 4106.5
$ws.Range("M31").Value = -3811.5

# Row 34
$ws.Range("H34").Value = 5550.9165
$ws.Range("I34").Value = 4106.5
$ws.Range("K34").Value = 4106.5
$ws.Range("M34").Value = -3904.5

# Row 62
$ws.Range("H62").Value = 4972.6
$ws.Range("I62").Value = 5987.6665
$ws.Range("J62").Value = 3450
$ws.Range("K62").Value = 5987.6665
$ws.Range("L62").Value = 3450
$ws.Range("M62").Value = -5363.6665
$ws.Range("N62").Value = -4698

# Row 65
$ws.Range("H65").Value = 4972.6
$ws.Range("I65").Value = 5987.6665
$ws.Range("J65").Value = 3450
$ws.Range("K65").Value = 29938.3325
$ws.Range("L65").Value = 17250
$ws.Range("M65").Value = -26818.3325
$ws.Range("N65").Value = -23490

# Row 132
$ws.Range("H132").Value = 16668327
$ws.Range("I132").Value = 18520200
$ws.Range("J132").Value = 1465.5
$ws.Range("K132").Value = 55560600
$ws.Range("L132").Value = 4396.5
$ws.Range("M132").Value = -55558070
$ws.Range("N132").Value = -9456.5

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 525
$ws.Range("I46").Value = 562.5
$ws.Range("K46").Value = 1687.5
$ws.Range("M46").Value = -1596.5

# Row 70
$ws.Range("H70").Value = 10133.333
$ws.Range("I70").Value = 5200
$ws.Range("K70").Value = 15600
$ws.Range("M70").Value = -15285

# Row 73
$ws.Range("H73").Value = 10133.333
$ws.Range("I73").Value = 5200
$ws.Range("K73").Value = 15600
$ws.Range("M73").Value = -14508

# Row 88
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856

# Row 91
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964

# Row 118
$ws.Range("H118").Value = 800
$ws.Range("I118").Value = 100
$ws.Range("J118").Value = 1500
$ws.Range("K118").Value = 300
$ws.Range("L118").Value = 4500
$ws.Range("M118").Value = 943
$ws.Range("N118").Value = -6986

# Row 122
$ws.Range("H122").Value = 1737.4
$ws.Range("I122").Value = 1009.1667
$ws.Range("K122").Value = 9082.5003
$ws.Range("M122").Value = -6632.5003

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1259.2858
$ws.Range("I97").Value = 1034.6666
$ws.Range("K97").Value = 1034.6666
$ws.Range("M97").Value = -538.6666

# Row 126
$ws.Range("H126").Value = 9135.817999999999
$ws.Range("I126").Value = 9552.143
$ws.Range("K126").Value = 28656.429
$ws.Range("M126").Value = -26186.429

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 20000
$ws.Range("J2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("N2").Value = -20224

# Row 136
$ws.Range("H136").Value = 2936.1428
$ws.Range("I136").Value = 2821
$ws.Range("J136").Value = 3000.111
$ws.Range("K136").Value = 8463
$ws.Range("L136").Value = 9000.332999999999
$ws.Range("M136").Value = -5913
$ws.Range("N136").Value = -14100.333

$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Range("H19").Value = 13497.5
$ws.Range("J19").Value = 13497.5
$ws.Range("L19").Value = 13497.5
$ws.Range("N19").Value = -13845.5

# Row 126
$ws.Range("H126").Value = 2394.375
$ws.Range("I126").Value = 2394.375
$ws.Range("K126").Value = 7183.125
$ws.Range("M126").Value = -4713.125

# Row 132
$ws.Range("H132").Value = 15158104
$ws.Range("I132").Value = 20003954
$ws.Range("K132").Value = 60011862
$ws.Range("M132").Value = -60009332
